# Update "想去人数" (want-to-go count) figures in column F across the four
# worksheets: 展览 (1), 演出 (2), 本地生活 (3), 全部类型 (4).
# Sheet index is used instead of sheet name to avoid any encoding issues.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$ws1.Range("F3").Value  = 981
$ws1.Range("F5").Value  = 439
$ws1.Range("F6").Value  = 682
$ws1.Range("F7").Value  = 242
$ws1.Range("F9").Value  = 17
$ws1.Range("F10").Value = 385
$ws1.Range("F11").Value = 186
$ws1.Range("F12").Value = 60
$ws1.Range("F13").Value = 782
$ws1.Range("F14").Value = 107
$ws1.Range("F15").Value = 1937
$ws1.Range("F16").Value = 437
$ws1.Range("F17").Value = 6302
$ws1.Range("F18").Value = 493
$ws1.Range("F19").Value = 514
$ws1.Range("F20").Value = 42
$ws1.Range("F21").Value = 82
$ws1.Range("F22").Value = 12
$ws1.Range("F23").Value = 197

# --- Sheet 2: 演出 ---
$ws2.Range("F4").Value = 29
$ws2.Range("F6").Value = 507

# --- Sheet 3: 本地生活 ---
$ws3.Range("F2").Value = 5443
$ws3.Range("F3").Value = 369
$ws3.Range("F4").Value = 367

# --- Sheet 4: 全部类型 ---
$ws4.Range("F3").Value  = 5443
$ws4.Range("F4").Value  = 369
$ws4.Range("F5").Value  = 367
$ws4.Range("F8").Value  = 29
$ws4.Range("F10").Value = 507
$ws4.Range("F11").Value = 507
$ws4.Range("F12").Value = 981
$ws4.Range("F16").Value = 439
$ws4.Range("F17").Value = 682
$ws4.Range("F18").Value = 242
$ws4.Range("F21").Value = 17
$ws4.Range("F22").Value = 385
$ws4.Range("F23").Value = 186
$ws4.Range("F25").Value = 60
$ws4.Range("F27").Value = 782
$ws4.Range("F28").Value = 107
$ws4.Range("F30").Value = 1937
$ws4.Range("F31").Value = 437
$ws4.Range("F32").Value = 6303
$ws4.Range("F34").Value = 493
$ws4.Range("F35").Value = 514
$ws4.Range("F36").Value = 42
$ws4.Range("F37").Value = 82
$ws4.Range("F39").Value = 12
$ws4.Range("F40").Value = 197
